$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2324.5417
$ws.Range("J40").Value = 2428.5
$ws.Range("L40").Value = 2428.5
$ws.Range("N40").Value = -2778.5
$ws.Range("H45").Value = 500
$ws.Range("I45").Value = 500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1308
$ws.Range("H53").Value = 120
$ws.Range("I53").Value = 120
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 120
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = 517
$ws.Range("H80").Value = 3320.6667
$ws.Range("I80").Value = 3320.6667
$ws.Range("K80").Value = 9962.000100000001
$ws.Range("M80").Value = -8964.000100000001
$ws.Range("H83").Value = 3320.6667
$ws.Range("I83").Value = 3320.6667
$ws.Range("K83").Value = 29886.0003
$ws.Range("M83").Value = -24894.0003
$ws.Range("H92").Value = 112051.445
$ws.Range("I92").Value = 112051.445
$ws.Range("K92").Value = 112051.445
$ws.Range("M92").Value = -110803.445
$ws.Range("H100").Value = 1675.7778
$ws.Range("I100").Value = 221.75
$ws.Range("J100").Value = 2839
$ws.Range("K100").Value = 221.75
$ws.Range("L100").Value = 2839
$ws.Range("M100").Value = 319.25
$ws.Range("N100").Value = -3921
$ws.Range("H132").Value = 1263.2222
$ws.Range("I132").Value = 1235.174
$ws.Range("J132").Value = 1424.5
$ws.Range("K132").Value = 3705.522
$ws.Range("L132").Value = 4273.5
$ws.Range("M132").Value = -1175.522
$ws.Range("N132").Value = -9333.5
$ws.Range("H137").Value = 1880.4
$ws.Range("J137").Value = 1880
$ws.Range("L137").Value = 5640
$ws.Range("N137").Value = -10740

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3514.875
$ws.Range("J2").Value = 4117.2
$ws.Range("L2").Value = 4117.2
$ws.Range("N2").Value = -4343.2
$ws.Range("H5").Value = 50
$ws.Range("I5").Value = 50
$ws.Range("K5").Value = 50
$ws.Range("M5").Value = 62
$ws.Range("H45").Value = 1993.4286
$ws.Range("I45").Value = 1993.4286
$ws.Range("K45").Value = 1993.4286
$ws.Range("M45").Value = -1616.4286
$ws.Range("H60").Value = 17398.6
$ws.Range("I60").Value = 14249.25
$ws.Range("J60").Value = 29996
$ws.Range("K60").Value = 14249.25
$ws.Range("L60").Value = 29996
$ws.Range("M60").Value = -13516.25
$ws.Range("N60").Value = -31462
$ws.Range("H116").Value = 3514.875
$ws.Range("J116").Value = 4117.2
$ws.Range("L116").Value = 4117.2
$ws.Range("N116").Value = -8705.200000000001
$ws.Range("H132").Value = 1224.75
$ws.Range("I132").Value = 1099.6897
$ws.Range("J132").Value = 1742.8572
$ws.Range("K132").Value = 3299.0691
$ws.Range("L132").Value = 5228.571599999999
$ws.Range("M132").Value = -769.0690999999997
$ws.Range("N132").Value = -10288.5716

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3514.875
$ws.Range("J3").Value = 4117.2
$ws.Range("L3").Value = 4117.2
$ws.Range("N3").Value = -4345.2
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("M4").Value = 65
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("H86").Value = 3032.4443
$ws.Range("I86").Value = 2538.9333
$ws.Range("K86").Value = 2538.9333
$ws.Range("M86").Value = -1415.9333
$ws.Range("H89").Value = 3032.4443
$ws.Range("I89").Value = 2538.9333
$ws.Range("K89").Value = 12694.6665
$ws.Range("M89").Value = -7078.666500000001
$ws.Range("H94").Value = 2116
$ws.Range("I94").Value = 2449
$ws.Range("K94").Value = 2449
$ws.Range("M94").Value = -1998

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 10449.5
$ws.Range("I25").Value = 900
$ws.Range("J25").Value = 19999
$ws.Range("K25").Value = 900
$ws.Range("L25").Value = 19999
$ws.Range("M25").Value = -726
$ws.Range("N25").Value = -20347
$ws.Range("H62").Value = 52819.875
$ws.Range("I62").Value = 3134.3333
$ws.Range("J62").Value = 82631.2
$ws.Range("K62").Value = 3134.3333
$ws.Range("L62").Value = 82631.2
$ws.Range("M62").Value = -2510.3333
$ws.Range("N62").Value = -83879.2
$ws.Range("H65").Value = 52819.875
$ws.Range("I65").Value = 3134.3333
$ws.Range("J65").Value = 82631.2
$ws.Range("K65").Value = 15671.6665
$ws.Range("L65").Value = 413156
$ws.Range("M65").Value = -12551.6665
$ws.Range("N65").Value = -419396
$ws.Range("H99").Value = 12654.863
$ws.Range("I99").Value = 8744.643
$ws.Range("J99").Value = 19497.75
$ws.Range("K99").Value = 8744.643
$ws.Range("L99").Value = 19497.75
$ws.Range("M99").Value = -7246.643
$ws.Range("N99").Value = -22493.75
$ws.Range("H107").Value = 656.4
$ws.Range("I107").Value = 618.2222
$ws.Range("K107").Value = 618.2222
$ws.Range("M107").Value = 1301.7778
$ws.Range("H122").Value = 3077
$ws.Range("I122").Value = 3756.7778
$ws.Range("J122").Value = 2312.25
$ws.Range("K122").Value = 11270.3334
$ws.Range("L122").Value = 6936.75
$ws.Range("M122").Value = -8820.3334
$ws.Range("N122").Value = -11836.75
$ws.Range("H126").Value = 12654.863
$ws.Range("I126").Value = 8744.643
$ws.Range("J126").Value = 19497.75
$ws.Range("K126").Value = 26233.929
$ws.Range("L126").Value = 58493.25
$ws.Range("M126").Value = -23763.929
$ws.Range("N126").Value = -63433.25
$ws.Range("H132").Value = 1487.2
$ws.Range("J132").Value = 1497.5
$ws.Range("L132").Value = 4492.5
$ws.Range("N132").Value = -9552.5
$ws.Range("H134").Value = 2090.5557
$ws.Range("I134").Value = 1545.2858
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 4635.857400000001
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -2100.857400000001
$ws.Range("N134").Value = -17067

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 179.5
$ws.Range("I60").Value = 179.5
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 538.5
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -287.5
$ws.Range("H113").Value = 3899.75
$ws.Range("I113").Value = 4899
$ws.Range("J113").Value = 3566.6667
$ws.Range("K113").Value = 14697
$ws.Range("L113").Value = 10700.0001
$ws.Range("M113").Value = -12527
$ws.Range("N113").Value = -15040.0001

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("N49").Value = 0
$ws.Range("H102").Value = 2619
$ws.Range("J102").Value = 2687.7778
$ws.Range("L102").Value = 2687.7778
$ws.Range("N102").Value = -5931.7778
$ws.Range("H109").Value = 26500
$ws.Range("J109").Value = 26500
$ws.Range("L109").Value = 26500
$ws.Range("N109").Value = -28580

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3432.0881
$ws.Range("I22").Value = 2127.276
$ws.Range("J22").Value = 11000
$ws.Range("K22").Value = 2127.276
$ws.Range("L22").Value = 11000
$ws.Range("M22").Value = -1832.276
$ws.Range("N22").Value = -11590
$ws.Range("H27").Value = 3432.0881
$ws.Range("I27").Value = 2127.276
$ws.Range("J27").Value = 11000
$ws.Range("K27").Value = 2127.276
$ws.Range("L27").Value = 11000
$ws.Range("M27").Value = -2020.276
$ws.Range("N27").Value = -11214
$ws.Range("H53").Value = 5000
$ws.Range("I53").Value = 5000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 5000
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -4482
$ws.Range("H82").Value = 949.8
$ws.Range("I82").Value = 874.5
$ws.Range("J82").Value = 1000
$ws.Range("K82").Value = 874.5
$ws.Range("L82").Value = 1000
$ws.Range("M82").Value = -513.5
$ws.Range("N82").Value = -1722
$ws.Range("H85").Value = 949.8
$ws.Range("I85").Value = 874.5
$ws.Range("J85").Value = 1000
$ws.Range("K85").Value = 874.5
$ws.Range("L85").Value = 1000
$ws.Range("M85").Value = 373.5
$ws.Range("N85").Value = -3496
$ws.Range("H88").Value = 26229.5
$ws.Range("J88").Value = 31582.334
$ws.Range("L88").Value = 31582.334
$ws.Range("N88").Value = -32438.334
$ws.Range("H91").Value = 26229.5
$ws.Range("J91").Value = 31582.334
$ws.Range("L91").Value = 31582.334
$ws.Range("N91").Value = -34546.334

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("N50").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("N53").Value = 0
$ws.Range("H81").Value = 7900.476
$ws.Range("I81").Value = 1536
$ws.Range("J81").Value = 13686.363
$ws.Range("K81").Value = 3072
$ws.Range("L81").Value = 27372.726
$ws.Range("M81").Value = -2011
$ws.Range("N81").Value = -29494.726
$ws.Range("H84").Value = 7900.476
$ws.Range("I84").Value = 1536
$ws.Range("J84").Value = 13686.363
$ws.Range("K84").Value = 15360
$ws.Range("L84").Value = 136863.63
$ws.Range("M84").Value = -10056
$ws.Range("N84").Value = -147471.63
$ws.Range("H120").Value = 39900
$ws.Range("J120").Value = 39900
$ws.Range("L120").Value = 39900
$ws.Range("N120").Value = -49576
$ws.Range("H132").Value = 23877.146
$ws.Range("I132").Value = 34052.645
$ws.Range("K132").Value = 102157.935
$ws.Range("M132").Value = -99627.935
